$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value to a cell while always keeping it as literal TEXT
# (mirrors how these "Price"/"Volume" columns are stored in the source sheet -
# plain text, even when the text happens to look like a number, e.g. "685.34"
# or "69.336.35"). Forcing NumberFormat to "@" before the assignment stops Excel
# from auto-coercing numeric-looking strings into real numbers; ClearFormats()
# afterwards restores the cell to the workbook default (General) formatting so
# no stray style is left behind on the cell.
function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

Set-TextValue "D2" "69.336.35"
Set-TextValue "E2" "  +0.25%  "
Set-TextValue "D3" "3.672.56"
Set-TextValue "E3" "  -0.21%  "
Set-TextValue "E4" "  +0.03%  "
Set-TextValue "D5" "685.34"
Set-TextValue "E5" "  +0.49%  "
Set-TextValue "D6" "159.19"
Set-TextValue "E6" "  -1.96%  "
Set-TextValue "E7" "  +0.07%  "
Set-TextValue "E8" "  -1.01%  "
Set-TextValue "E9" "  -1.57%  "
Set-TextValue "D10" "7.08"
Set-TextValue "E10" "  -2.19%  "
Set-TextValue "E11" "  -3.59%  "
Set-TextValue "E12" "  -1.14%  "
Set-TextValue "D13" "4.294.61"
Set-TextValue "E13" "  -0.16%  "
Set-TextValue "D14" "32.14"
Set-TextValue "E14" "  -3.83%  "
Set-TextValue "D15" "3.673.54"
Set-TextValue "E15" "  -0.06%  "
Set-TextValue "D16" "69.313.81"
Set-TextValue "E16" "  +0.15%  "
Set-TextValue "E17" "  +2.04%  "
Set-TextValue "D18" "15.78"
Set-TextValue "E18" "  -3.27%  "
Set-TextValue "E19" "  -4.00%  "
Set-TextValue "D20" "469.63"
Set-TextValue "E20" "  -2.74%  "
Set-TextValue "E21" "  +0.93%  "
Set-TextValue "E22" "  -2.48%  "
Set-TextValue "D23" "79.54"
Set-TextValue "E23" "  -0.20%  "
Set-TextValue "D24" "3.822.63"
Set-TextValue "E24" "  -0.04%  "
Set-TextValue "E25" "  +0.08%  "
Set-TextValue "E26" "  -2.94%  "
Set-TextValue "D27" "10.94"
Set-TextValue "E27" "  -5.22%  "
Set-TextValue "D28" "9.18"
Set-TextValue "E28" "  -3.73%  "
Set-TextValue "E29" "  -1.34%  "
Set-TextValue "E30" "  -5.64%  "
Set-TextValue "D31" "6.56"
Set-TextValue "E31" "  -2.63%  "
Set-TextValue "B32" "Binance-PegBSC-USD"
Set-TextValue "C32" "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextValue "D32" "0.999"
Set-TextValue "E32" "  -0.04%  "
Set-TextValue "B33" "ImmutableX"
Set-TextValue "C33" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D33" "1.99"
Set-TextValue "E33" "  -5.80%  "
Set-TextValue "D34" "26.81"
Set-TextValue "E34" "  -0.49%  "
Set-TextValue "D35" "3.647.67"
Set-TextValue "E35" "  +0.02%  "
Set-TextValue "D36" "0.159"
Set-TextValue "E36" "  -1.75%  "
Set-TextValue "D37" "8.15"
Set-TextValue "E37" "  -4.20%  "
Set-TextValue "E38" "  +0.88%  "
Set-TextValue "D40" "2.21"
Set-TextValue "E40" "  +1.70%  "
Set-TextValue "D41" "0.0895"
Set-TextValue "E41" "  -5.34%  "
Set-TextValue "E42" "  +0.01%  "
Set-TextValue "E43" "  -1.73%  "
Set-TextValue "D44" "165.83"
Set-TextValue "E44" "  +5.76%  "
Set-TextValue "D45" "47.49"
Set-TextValue "E45" "  -1.01%  "
Set-TextValue "D46" "0.000280"
Set-TextValue "E46" "  +0.99%  "
Set-TextValue "E47" "  -2.75%  "
Set-TextValue "D48" "1.11"
Set-TextValue "E48" "  +5.13%  "
Set-TextValue "D49" "1.29"
Set-TextValue "E49" "  -0.85%  "
Set-TextValue "D50" "27.65"
Set-TextValue "E50" "  -1.77%  "
Set-TextValue "D51" "7.75"
Set-TextValue "E51" "  -3.59%  "
